$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: "Means" ----------
$ws1 = $wb.Worksheets.Item("Means")

# New header columns F and G
$ws1.Range("F1").Value = "Within 5 miles of HFC production facility"
$ws1.Range("G1").Value = "Within 10 miles of HFC production facility"

# New column values for existing rows
$ws1.Range("F2").Value = 58
$ws1.Range("G2").Value = 62

$ws1.Range("F3").Value = 40
$ws1.Range("G3").Value = 31

$ws1.Range("F4").Value = 2.2
$ws1.Range("G4").Value = 7.4

$ws1.Range("F5").Value = 2.5
$ws1.Range("G5").Value = 5.2

$ws1.Range("F6").Value = 78
$ws1.Range("G6").Value = 82

$ws1.Range("F7").Value = 2.8
$ws1.Range("G7").Value = 6.2

$ws1.Range("F8").Value = 4.6
$ws1.Range("G8").Value = 5.3

# Row 9 (Total Cancer Risk) updated existing values + new columns
$ws1.Range("B9").Value = 26
$ws1.Range("C9").Value = 39
$ws1.Range("D9").Value = 200
$ws1.Range("E9").Value = 130
$ws1.Range("F9").Value = 120
$ws1.Range("G9").Value = 82

# Row 10 (Total Respiratory) updated existing values + new columns
$ws1.Range("B10").Value = 0.32
$ws1.Range("C10").Value = 0.43
$ws1.Range("D10").Value = 0.6
$ws1.Range("E10").Value = 0.55
$ws1.Range("F10").Value = 0.54
$ws1.Range("G10").Value = 0.54

# ---------- Sheet 2: "Standard Deviations" ----------
$ws2 = $wb.Worksheets.Item("Standard Deviations")

# New header columns F and G
$ws2.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$ws2.Range("G1").Value = "Within 10 mile of HFC production facility SD"

# New column values for existing rows
$ws2.Range("F2").Value = 24
$ws2.Range("G2").Value = 26

$ws2.Range("F3").Value = 24
$ws2.Range("G3").Value = 27

$ws2.Range("F4").Value = 2.3
$ws2.Range("G4").Value = 6.5

$ws2.Range("F5").Value = 9.7
$ws2.Range("G5").Value = 6.3

$ws2.Range("F6").Value = 35
$ws2.Range("G6").Value = 36

$ws2.Range("F7").Value = 2.6
$ws2.Range("G7").Value = 8.1

$ws2.Range("F8").Value = 10
$ws2.Range("G8").Value = 7.8

# Row 9 (Total Cancer Risk) updated existing values + new columns
$ws2.Range("B9").Value = 8.6
$ws2.Range("C9").Value = 24
$ws2.Range("E9").Value = 64
$ws2.Range("F9").Value = 62
$ws2.Range("G9").Value = 32

# Row 10 (Total Respiratory) updated existing values + new columns
$ws2.Range("B10").Value = 0.14
$ws2.Range("C10").Value = 0.084
$ws2.Range("E10").Value = 0.095
$ws2.Range("F10").Value = 0.092
$ws2.Range("G10").Value = 0.08
